$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.84
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 1.44
$ws.Range("N2").Value = 3.3
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.79
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.3
$ws.Range("S2").Value = 3.85
$ws.Range("T2").Value = 1.95
$ws.Range("U2").Value = 1.9
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 17
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 130
$ws.Range("AB2").Value = 8
$ws.Range("AC2").Value = 8.2
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 10.5
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 22
$ws.Range("AJ2").Value = 20
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 42
$ws.Range("AM2").Value = 140
$ws.Range("AN2").Value = 15
$ws.Range("AO2").Value = 95

# Row 3
$ws.Range("P3").Value = 1.76
$ws.Range("Q3").Value = 2.28
$ws.Range("R3").Value = 1.29
$ws.Range("S3").Value = 4.1
$ws.Range("T3").Value = 2.02
$ws.Range("U3").Value = 1.88
$ws.Range("Y3").Value = 7.6
$ws.Range("AA3").Value = 1000
$ws.Range("AE3").Value = 29
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AO3").Value = 17

# Row 4
$ws.Range("N4").Value = 2.94
$ws.Range("AG4").Value = 24

# Row 5
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 8.6
$ws.Range("AI5").Value = 120

# Row 6
$ws.Range("K6").Value = 3.7

# Row 7
$ws.Range("F7").Value = 3.25
$ws.Range("H7").Value = 2.4
$ws.Range("R7").Value = 1.36
$ws.Range("T7").Value = 1.76
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 38
$ws.Range("AN7").Value = 1000

# Row 8
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 3.95
$ws.Range("N8").Value = 3.65
$ws.Range("P8").Value = 2.12
$ws.Range("S8").Value = 2.78
$ws.Range("T8").Value = 1.76
$ws.Range("U8").Value = 2.1
$ws.Range("X8").Value = 20
$ws.Range("AB8").Value = 12
$ws.Range("AF8").Value = 11.5
$ws.Range("AG8").Value = 12.5
$ws.Range("AI8").Value = 75
$ws.Range("AJ8").Value = 18
$ws.Range("AK8").Value = 17.5
$ws.Range("AN8").Value = 9.2
$ws.Range("AO8").Value = 85
